$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.472.64"

$ws.Range("D3").Value = "'1.571.39"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "'288.56"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").Value = "'0.3727"
$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("D8").Value = "'48.37"
$ws.Range("E8").Value = "  -2.86%  "

$ws.Range("D9").Value = "'0.3333"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").Value = "'0.07487"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("D11").Value = "'1.132"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "'20.91"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").Value = "'5.974"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "'6.917"
$ws.Range("E15").Value = "  -0.67%  "

$ws.Range("D16").Value = "'1.576.97"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").Value = "'0.00001117"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "'88.29"
$ws.Range("E18").Value = "  -2.29%  "

$ws.Range("D19").Value = "'0.06760"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").Value = "'16.46"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("D23").Value = "'12.09"
$ws.Range("E23").Value = "  -0.93%  "

$ws.Range("D24").Value = "'22.472.33"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "'2.388"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").Value = "'2.569"
$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("D27").Value = "'152.69"
$ws.Range("E27").Value = "  +2.50%  "

$ws.Range("D28").Value = "'19.73"
$ws.Range("E28").Value = "  -1.23%  "

$ws.Range("D29").Value = "'5.006"
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").Value = "'124.01"
$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").Value = "'1.752.32"
$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("D32").Value = "'1.053"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("D33").Value = "'6.158"
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("D35").Value = "'9.676"
$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("D36").Value = "'0.08305"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").Value = "'0.02459"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").Value = "'0.2272"
$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("D39").Value = "'0.06377"
$ws.Range("E39").Value = "  -2.71%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.295"
$ws.Range("E40").Value = "  -4.57%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "'5.376"
$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("D42").Value = "'11.30"
$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("D43").Value = "'0.6304"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Value = "'13.90"
$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("D46").Value = "'0.6141"
$ws.Range("E46").Value = "  +5.10%  "

$ws.Range("D47").Value = "'3.782"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").Value = "'2.049"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("D49").Value = "'125.07"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").Value = "'1.214"
$ws.Range("E50").Value = "  -1.98%  "

$ws.Range("D51").Value = "'0.07265"
$ws.Range("E51").Value = "  -0.47%  "
